{"js": "// Apply the two content edits described in the commit:\n//\n// 1) In the \"-There is a field in the database...\" paragraph, the three\n//    runs split around the grammar-checker markers (\"look into\") are\n//    merged back into one run with the same combined text (no wording\n//    change, just de-fragmenting the runs / dropping the proofErr marks).\n//\n// 2) In the red \"I need to write a prototype...\" paragraph, the runs\n//    around \"at the moment\" are likewise merged into the preceding\n//    \"Write a function...\" run, and the paragraph is then split right\n//    after that merged sentence, with a brand-new paragraph (same red\n//    formatting, no strike-through) added containing the new sentence\n//    about query strings / the forward button.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1: de-fragment the runs in the \"-There is a field...\" paragraph.\n// ---------------------------------------------------------------------\nconst search1 = body.search(\n  \"-There is a field in the database for the key part of the question.\",\n  { matchCase: true }\n);\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  const para1 = search1.items[0].paragraphs.getFirst();\n  const para1Range = para1.getRange();\n  para1Range.insertText(\n    \"-There is a field in the database for the key part of the question. \" +\n      \"I could use this to keep the meaning (and answer) of the question \" +\n      \"the same while changing the wording of it? I will have to look \" +\n      \"into methods of doing this if it is possible.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 2: merge the \"Write a function...\"/\"at the moment\"/\"...system.\"\n// runs into a single run, then split the paragraph right after it and\n// add the new \"Because the website...\" paragraph.\n// ---------------------------------------------------------------------\nconst search2 = body.search(\"Write a function to fetch the question, but\", {\n  matchCase: true,\n});\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  const found = search2.items[0];\n  const hostParagraph = found.paragraphs.getFirst();\n\n  // Expand the found range so it covers everything from \"Write a\n  // function...\" through to the end of the paragraph (i.e. also the\n  // \"at the moment\" / \"...rewording system.\" runs + the proofErr marks\n  // sitting between them).\n  const paragraphEnd = hostParagraph.getRange(Word.RangeLocation.end);\n  const toMerge = found.expandTo(paragraphEnd);\n\n  // Insert the fully-merged replacement sentence right after the range\n  // we are about to delete, then delete the original fragmented runs.\n  // Inserting first (instead of using InsertLocation.replace) keeps the\n  // new text out of the about-to-be-removed proofErr-bounded runs, so\n  // the stray empty <w:proofErr/> markers get cleaned up along with the\n  // rest of the deleted range instead of being left behind.\n  const insertionPoint = toMerge.getRange(Word.RangeLocation.end);\n  insertionPoint.insertText(\n    \"Write a function to fetch the question, but at the moment the \" +\n      \"function can just return what it is inputted just to make it \" +\n      \"easier when I implement the question rewording system.\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n\n  toMerge.delete();\n  await context.sync();\n\n  // Now split the paragraph: add a brand-new paragraph right after the\n  // host paragraph with the new sentence (inherits the red, non-strike\n  // character formatting from the insertion point).\n  const hostEnd = hostParagraph.getRange(Word.RangeLocation.end);\n  hostEnd.insertParagraph(\n    \"Because the website uses query strings to pass quiz ID and such it \" +\n      \"does break if the forward button is used potentially\u2026.. This is \" +\n      \"something that could be fixed but is not my priority\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "# Apply the two content edits described in the commit:\n#\n# 1) In the \"-There is a field in the database...\" paragraph, the three\n#    runs split around the grammar-checker markers (\"look into\") are\n#    merged back into one run with the same combined text (no wording\n#    change, just de-fragmenting the runs / dropping the proofErr marks).\n#\n# 2) In the red \"I need to write a prototype...\" paragraph, the runs\n#    around \"at the moment\" are likewise merged into the preceding\n#    \"Write a function...\" run, and the paragraph is then split right\n#    after that merged sentence, with a brand-new paragraph (same red\n#    formatting, no strike-through) added containing the new sentence\n#    about query strings / the forward button.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: de-fragment the runs in the \"-There is a field...\" paragraph.\n# Find-and-replace with the identical (fully concatenated) text merges\n# the runs that were split by the gramStart/gramEnd proof-error markers\n# into a single run and drops those markers.\n# ---------------------------------------------------------------------\n$text1 = \"-There is a field in the database for the key part of the question. I could use this to keep the meaning (and answer) of the question the same while changing the wording of it? I will have to look into methods of doing this if it is possible.\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$found1 = $rng1.Find.Execute($text1, $false, $false, $false, $false, $false, $true, 1, $false, $text1, 2)\n\n# ---------------------------------------------------------------------\n# Change 2a: merge the \"Write a function...\" / \"at the moment\" /\n# \"...system.\" runs into a single run the same way.\n# ---------------------------------------------------------------------\n$text2 = \"Write a function to fetch the question, but at the moment the function can just return what it is inputted just to make it easier when I implement the question rewording system.\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute($text2, $false, $false, $false, $false, $false, $true, 1, $false, $text2, 2)\n\n# ---------------------------------------------------------------------\n# Change 2b: split the host paragraph right after that merged sentence\n# and add a brand-new paragraph with the new sentence about query\n# strings / the forward button (inherits the red, non-strike formatting\n# already on the paragraph mark).\n# ---------------------------------------------------------------------\n# $rng2 is now collapsed to the (re-merged) found text; locate which\n# paragraph it lives in by comparing character offsets rather than\n# relying on the Range's own .Paragraphs collection.\n$targetStart = $rng2.Start\n$hostIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($targetStart -ge $p.Range.Start -and $targetStart -lt $p.Range.End) {\n        $hostIndex = $i\n        break\n    }\n}\n\nif ($hostIndex -gt 0) {\n    $hostParagraph = $d.Paragraphs.Item($hostIndex)\n    $hostParagraph.Range.InsertParagraphAfter()\n    $newParagraph = $d.Paragraphs.Item($hostIndex + 1)\n    $newParagraph.Range.Text = \"Because the website uses query strings to pass quiz ID and such it does break if the forward button is used potentially\u2026.. This is something that could be fixed but is not my priority\"\n}\n\nWrite-Output \"done\"\n"}
